# Scotland Premiership - base update (19-06-2024 21:51)
# The underlying source rows were re-ordered/re-matched; apply that by
# relocating each affected row's match data (columns B:AD) to the row that
# now carries it, while leaving column A (the running index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

function Get-RowData($row) {
    $data = @{}
    foreach ($col in $cols) {
        $data[$col] = $ws.Range("$col$row").Value2
    }
    return $data
}

function Set-RowData($row, $data) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $data[$col]
    }
}

# target_row -> source_row (source row's B:AD content now belongs on target_row)
$map = @{
    107 = 108; 108 = 107;
    115 = 117; 116 = 118; 117 = 116; 118 = 115;
    119 = 120; 120 = 121; 121 = 119;
    136 = 137; 137 = 138; 138 = 136;
    143 = 144; 144 = 143;
    207 = 209; 208 = 207; 209 = 208;
    228 = 229; 229 = 228;
}

# Snapshot every affected row's current data before any writes, so a row
# used as a source later isn't clobbered first.
$snapshot = @{}
foreach ($row in $map.Keys) {
    $snapshot[$row] = Get-RowData $row
}

foreach ($row in $map.Keys) {
    $srcRow = $map[$row]
    Set-RowData $row $snapshot[$srcRow]
}
